$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "theta_threshold_range" (row 5) is removed entirely; this shifts
# the "pie_threshold_range" row (old row 6) up to become the new row 5.
$ws.Rows("5:5").Delete() | Out-Null

# Updated Min/Max values for the remaining parameter rows.
$ws.Range("B2").Value = 5.2
$ws.Range("C2").Value = 10.1
$ws.Range("B3").Value = 3.8
$ws.Range("C3").Value = 8.9
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.6
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# C4 and B5 (old B6) previously carried a one-off Times New Roman font; bring
# them back in line with the rest of the body cells by copying that format.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column widths: A narrower, B/C sized to fit their (now short) numeric content.
$ws.Columns("A").ColumnWidth = 20.714285714285715
$ws.Columns("B").ColumnWidth = 4.428571428571429
$ws.Columns("C").ColumnWidth = 4.857142857142857

# Selection moves to C3.
$ws.Range("C3").Select() | Out-Null
